# Lentils, faba beans and chickpeas share their stubble-simulation inputs with
# lupins. Those three landuses previously carried their own (now superseded)
# simulated stubble figures in columns G:R (first replicate block) and V:AD
# (second replicate block); the "lupins" results live in columns AK:AM and are
# left untouched. No profit-driver change results because lentils/faba/chickpeas
# are not selected land uses in this scenario -- only the raw simulation inputs
# move into line with lupins.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("G1").Value = 0.00614513094148918
$ws.Range("H1").Value = 0.1599499218526465
$ws.Range("I1").Value = 0.009544565079334256
$ws.Range("J1").Value = 0.09505056148880539
$ws.Range("K1").Value = 0.06354486258187295
$ws.Range("L1").Value = 0.145204900316941
$ws.Range("M1").Value = 0.00614513094148918
$ws.Range("N1").Value = 0.1599499218526465
$ws.Range("O1").Value = 0.009544565079334256
$ws.Range("P1").Value = 0.00614513094148918
$ws.Range("Q1").Value = 0.1599499218526465
$ws.Range("R1").Value = 0.009544565079334256
$ws.Range("V1").Value = 0.1458254340789792
$ws.Range("X1").Value = 0.1936065337559213
$ws.Range("Y1").Value = 0.1458254340789792
$ws.Range("AA1").Value = 0.1936065337559213
$ws.Range("AB1").Value = 0.1458213408483589
$ws.Range("AD1").Value = 0.1936010993390977

# Row 2
$ws.Range("G2").Value = 0.0005586482674081073
$ws.Range("H2").Value = 0.01454090198660423
$ws.Range("I2").Value = 0.0008676877344849324
$ws.Range("J2").Value = 0.01945401724727308
$ws.Range("K2").Value = 0.005776805689261176
$ws.Range("L2").Value = 0.02971911570966397
$ws.Range("M2").Value = 0.0005586482674081073
$ws.Range("N2").Value = 0.01454090198660423
$ws.Range("O2").Value = 0.0008676877344849324
$ws.Range("P2").Value = 0.0005586482674081073
$ws.Range("Q2").Value = 0.01454090198660423
$ws.Range("R2").Value = 0.0008676877344849324
$ws.Range("V2").Value = 0.029846120477654
$ws.Range("X2").Value = 0.0396254876128853
$ws.Range("Y2").Value = 0.029846120477654
$ws.Range("AA2").Value = 0.0396254876128853
$ws.Range("AB2").Value = 0.02984528271533216
$ws.Range("AD2").Value = 0.03962437534971758

# Row 3
$ws.Range("G3").Value = 0.1325050224731904
$ws.Range("H3").Value = 0.01454090198660423
$ws.Range("I3").Value = 0.2058056732030404
$ws.Range("J3").Value = 0.009623965327339321
$ws.Range("K3").Value = 0.1491236110031118
$ws.Range("L3").Value = 0.01470214277665879
$ws.Range("M3").Value = 0.1325050224731904
$ws.Range("N3").Value = 0.01454090198660423
$ws.Range("O3").Value = 0.2058056732030404
$ws.Range("P3").Value = 0.1325050224731904
$ws.Range("Q3").Value = 0.01454090198660423
$ws.Range("R3").Value = 0.2058056732030404
$ws.Range("V3").Value = 0.01476497244664452
$ws.Range("X3").Value = 0.01960285703554505
$ws.Range("Y3").Value = 0.01476497244664452
$ws.Range("AA3").Value = 0.01960285703554505
$ws.Range("AB3").Value = 0.0147645580029111
$ws.Range("AD3").Value = 0.01960230679535431

# Row 4
$ws.Range("G4").Value = 0.07240135459503878
$ws.Range("H4").Value = 0.0647871303291303
$ws.Range("I4").Value = 0.112453167775273
$ws.Range("J4").Value = 0.009623965327339321
$ws.Range("K4").Value = 0.1475814541468215
$ws.Range("L4").Value = 0.01470214277665879
$ws.Range("M4").Value = 0.07240135459503878
$ws.Range("N4").Value = 0.0647871303291303
$ws.Range("O4").Value = 0.112453167775273
$ws.Range("P4").Value = 0.07240135459503878
$ws.Range("Q4").Value = 0.0647871303291303
$ws.Range("R4").Value = 0.112453167775273
$ws.Range("V4").Value = 0.01476497244664452
$ws.Range("X4").Value = 0.01960285703554505
$ws.Range("Y4").Value = 0.01476497244664452
$ws.Range("AA4").Value = 0.01960285703554505
$ws.Range("AB4").Value = 0.0147645580029111
$ws.Range("AD4").Value = 0.01960230679535431

# Row 5
$ws.Range("G5").Value = 0.01799446241337307
$ws.Range("H5").Value = 0.2483404958213521
$ws.Range("I5").Value = 0.02794884587609009
$ws.Range("J5").Value = 0.009623965327339321
$ws.Range("K5").Value = 0.06060401653451566
$ws.Range("L5").Value = 0.01470214277665879
$ws.Range("M5").Value = 0.01799446241337307
$ws.Range("N5").Value = 0.2483404958213521
$ws.Range("O5").Value = 0.02794884587609009
$ws.Range("P5").Value = 0.01799446241337307
$ws.Range("Q5").Value = 0.2483404958213521
$ws.Range("R5").Value = 0.02794884587609009
$ws.Range("V5").Value = 0.01476497244664452
$ws.Range("X5").Value = 0.01960285703554505
$ws.Range("Y5").Value = 0.01476497244664452
$ws.Range("AA5").Value = 0.01960285703554505
$ws.Range("AB5").Value = 0.0147645580029111
$ws.Range("AD5").Value = 0.01960230679535431

# Row 6
$ws.Range("G6").Value = 0.01799446241337307
$ws.Range("H6").Value = 0.1603448127663619
$ws.Range("I6").Value = 0.02794884587609009
$ws.Range("J6").Value = 0.009623965327339321
$ws.Range("K6").Value = 0.03325034766821872
$ws.Range("L6").Value = 0.01470214277665879
$ws.Range("M6").Value = 0.01799446241337307
$ws.Range("N6").Value = 0.1603448127663619
$ws.Range("O6").Value = 0.02794884587609009
$ws.Range("P6").Value = 0.01799446241337307
$ws.Range("Q6").Value = 0.1603448127663619
$ws.Range("R6").Value = 0.02794884587609009
$ws.Range("V6").Value = 0.01476497244664452
$ws.Range("X6").Value = 0.01960285703554505
$ws.Range("Y6").Value = 0.01476497244664452
$ws.Range("AA6").Value = 0.01960285703554505
$ws.Range("AB6").Value = 0.0147645580029111
$ws.Range("AD6").Value = 0.01960230679535431

# Row 7
$ws.Range("G7").Value = 0.01799446241337307
$ws.Range("H7").Value = 0.0508931569531148
$ws.Range("I7").Value = 0.02794884587609009
$ws.Range("J7").Value = 0.009623965327339321
$ws.Range("K7").Value = 0.03325034766821872
$ws.Range("L7").Value = 0.01470214277665879
$ws.Range("M7").Value = 0.01799446241337307
$ws.Range("N7").Value = 0.0508931569531148
$ws.Range("O7").Value = 0.02794884587609009
$ws.Range("P7").Value = 0.01799446241337307
$ws.Range("Q7").Value = 0.0508931569531148
$ws.Range("R7").Value = 0.02794884587609009
$ws.Range("V7").Value = 0.01476497244664452
$ws.Range("X7").Value = 0.01960285703554505
$ws.Range("Y7").Value = 0.01476497244664452
$ws.Range("AA7").Value = 0.01960285703554505
$ws.Range("AB7").Value = 0.0147645580029111
$ws.Range("AD7").Value = 0.01960230679535431

# Row 8
$ws.Range("G8").Value = 0.01799446241337307
$ws.Range("H8").Value = 0.0508931569531148
$ws.Range("I8").Value = 0.02794884587609009
$ws.Range("J8").Value = 0.009623965327339321
$ws.Range("K8").Value = 0.03325034766821872
$ws.Range("L8").Value = 0.01470214277665879
$ws.Range("M8").Value = 0.01799446241337307
$ws.Range("N8").Value = 0.0508931569531148
$ws.Range("O8").Value = 0.02794884587609009
$ws.Range("P8").Value = 0.01799446241337307
$ws.Range("Q8").Value = 0.0508931569531148
$ws.Range("R8").Value = 0.02794884587609009
$ws.Range("V8").Value = 0.01476497244664452
$ws.Range("X8").Value = 0.01960285703554505
$ws.Range("Y8").Value = 0.01476497244664452
$ws.Range("AA8").Value = 0.01960285703554505
$ws.Range("AB8").Value = 0.0147645580029111
$ws.Range("AD8").Value = 0.01960230679535431

# Row 9
$ws.Range("G9").Value = 0.01799446241337307
$ws.Range("H9").Value = 0.0508931569531148
$ws.Range("I9").Value = 0.02794884587609009
$ws.Range("J9").Value = 0.009623965327339321
$ws.Range("K9").Value = 0.03325034766821872
$ws.Range("L9").Value = 0.01470214277665879
$ws.Range("M9").Value = 0.01799446241337307
$ws.Range("N9").Value = 0.0508931569531148
$ws.Range("O9").Value = 0.02794884587609009
$ws.Range("P9").Value = 0.01799446241337307
$ws.Range("Q9").Value = 0.0508931569531148
$ws.Range("R9").Value = 0.02794884587609009
$ws.Range("V9").Value = 0.01476497244664452
$ws.Range("X9").Value = 0.01960285703554505
$ws.Range("Y9").Value = 0.01476497244664452
$ws.Range("AA9").Value = 0.01960285703554505
$ws.Range("AB9").Value = 0.0147645580029111
$ws.Range("AD9").Value = 0.01960230679535431

# Row 10
$ws.Range("G10").Value = 0.02484226222688343
$ws.Range("H10").Value = 0.0508931569531148
$ws.Range("I10").Value = 0.03858479026728702
$ws.Range("J10").Value = 0.01329732916569063
$ws.Range("K10").Value = 0.09058906979375894
$ws.Range("L10").Value = 0.020313792214821
$ws.Range("M10").Value = 0.02484226222688343
$ws.Range("N10").Value = 0.0508931569531148
$ws.Range("O10").Value = 0.03858479026728702
$ws.Range("P10").Value = 0.02484226222688343
$ws.Range("Q10").Value = 0.0508931569531148
$ws.Range("R10").Value = 0.03858479026728702
$ws.Range("V10").Value = 0.02040060329266212
$ws.Range("X10").Value = 0.027085056286428
$ws.Range("Y10").Value = 0.02040060329266212
$ws.Range("AA10").Value = 0.027085056286428
$ws.Range("AB10").Value = 0.02040003066022254
$ws.Range("AD10").Value = 0.02708429602548693
